$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings need to be
# forced to text so Excel does not silently convert them to numbers
# (which would also reformat the displayed text, e.g. "0.0920" -> 0.092).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "44.142.16"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "2.349.84"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "313.94"
$ws.Range("E5").Value = "  +0.46%  "
Set-TextValue "D6" "109.45"
$ws.Range("E6").Value = "  +7.19%  "
Set-TextValue "D7" "0.629"
$ws.Range("E7").Value = "  +1.13%  "
Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +4.21%  "
Set-TextValue "D10" "41.86"
$ws.Range("E10").Value = "  +8.19%  "
Set-TextValue "D11" "0.0920"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("E13").Value = "  +4.39%  "
$ws.Range("E14").Value = "  -0.33%  "
Set-TextValue "D15" "15.53"
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("D16").Value = "2.703.03"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "2.345.51"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "44.117.23"
$ws.Range("E18").Value = "  +4.04%  "
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("E20").Value = "  +2.57%  "
Set-TextValue "D21" "12.96"
$ws.Range("E21").Value = "  -2.65%  "
Set-TextValue "D22" "74.71"
$ws.Range("E22").Value = "  +2.54%  "
Set-TextValue "D23" "3.51"
$ws.Range("E23").Value = "  +0.87%  "
Set-TextValue "D24" "266.81"
$ws.Range("E24").Value = "  +1.38%  "
Set-TextValue "D25" "2.28"
$ws.Range("E25").Value = "  +5.73%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.78%  "
Set-TextValue "D27" "7.60"
$ws.Range("E27").Value = "  +9.55%  "
Set-TextValue "D28" "11.20"
$ws.Range("E28").Value = "  +5.25%  "
$ws.Range("E29").Value = "  +2.63%  "
Set-TextValue "D30" "39.59"
$ws.Range("E30").Value = "  +10.54%  "
$ws.Range("E31").Value = "  +0.87%  "
Set-TextValue "D32" "169.19"
$ws.Range("E32").Value = "  +2.59%  "
Set-TextValue "D33" "0.0915"
$ws.Range("E33").Value = "  +6.15%  "
Set-TextValue "D34" "2.81"
$ws.Range("E34").Value = "  +7.80%  "
$ws.Range("E35").Value = "  +1.32%  "
Set-TextValue "D36" "0.117"
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("E37").Value = "  +6.07%  "
$ws.Range("E38").Value = "  +5.55%  "
$ws.Range("E39").Value = "  +10.27%  "
$ws.Range("E40").Value = "  +3.34%  "
Set-TextValue "D41" "1.74"
$ws.Range("E41").Value = "  +10.77%  "
Set-TextValue "D42" "104.16"
$ws.Range("E42").Value = "  +5.57%  "
Set-TextValue "D43" "13.98"
$ws.Range("E43").Value = "  +17.54%  "
$ws.Range("E44").Value = "  +5.91%  "
Set-TextValue "D45" "71.18"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("E46").Value = "  -0.53%  "
Set-TextValue "D47" "116.12"
$ws.Range("E47").Value = "  +5.97%  "
Set-TextValue "D48" "78.12"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "9.04"
$ws.Range("E49").Value = "  +4.83%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.653.90"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D51" "5.35"
$ws.Range("E51").Value = "  +3.28%  "
